$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2023-01-02 12:57:01"

# Column O = 15 (timestamp), update every data row (2..397) to new timestamp
$lastRow = 397
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 15).Value = $newTimestamp
}

# Column M = 13 (productAriaLabel), append promo text to specific rows
$mUpdates = @{
    25  = "Dar-Vida Cracker Thymian &amp; Salz 2+1 Aktion 3.95 Schweizer Franken"
    29  = "Dar-Vida Cracker Käse 2+1 Aktion 3.95 Schweizer Franken"
    32  = "Roland Zwieback Classic 20% ab 2 Aktion 4.40 Schweizer Franken"
    39  = "Ölz Premium Drei-Korn Toast - Online kein Bestand 3.40 Schweizer Franken"
    105 = "Betty Bossi Naturaplan Bio Kuchenteig rund ausgewallt Ø32cm - Online kein Bestand 2.10 Schweizer Franken"
    107 = "Dar-Vida Cracker Tomate &amp; Basilikum 2+1 Aktion 3.95 Schweizer Franken"
    133 = "Dar-Vida Break Cranberry &amp; Apfel 2+1 Aktion 3.10 Schweizer Franken"
    136 = "Dar-Vida Break Choco &amp; Cacaonibs 2+1 Aktion 3.10 Schweizer Franken"
    139 = "Dar-Vida Cracker Schokolade dunkel 2+1 Aktion 4.95 Schweizer Franken"
    141 = "Dar-Vida Cracker Ur-Dinkel 2+1 Aktion 3.95 Schweizer Franken"
    151 = "Dar-Vida Honig extra fin 2+1 Aktion 3.95 Schweizer Franken"
    157 = "Dar-Vida Sandwich Frischkäse 2+1 Aktion 5.20 Schweizer Franken"
    166 = "Dar-Vida Oliven extra fin 2+1 Aktion 3.95 Schweizer Franken"
    186 = "Dar-Vida Naturaplan Bio Cracker Kürbiskerne 2+1 Aktion 3.95 Schweizer Franken"
    190 = "Dar-Vida Cracker Nature 2+1 Aktion 3.60 Schweizer Franken"
    193 = "Roland Knäckebrot Hafer 20% ab 2 Aktion 3.95 Schweizer Franken"
    200 = "Dar-Vida Sandwich Schokolade &amp; Haselnusscrème 2+1 Aktion 5.20 Schweizer Franken"
    203 = "Bonne Maman Madeleine au Citron 7 Stück 20% ab 2 Aktion 3.95 Schweizer Franken"
    237 = "Dar-Vida Sandwich Tomate &amp; Basilicum 2+1 Aktion 5.20 Schweizer Franken"
    240 = "Dar-Vida Cracker Chia &amp; Quinoa 2+1 Aktion 3.95 Schweizer Franken"
}

foreach ($row in $mUpdates.Keys) {
    $ws.Cells.Item($row, 13).Value = $mUpdates[$row]
}

Write-Host "Done"
